{"js": "// Update the multiplication problems/answers table: each cell's old\n// \"A\u00d7B=C\" text is replaced with a new \"A\u00d7B=C\" text, preserving the\n// run's existing formatting (font, size, etc.) by replacing only the\n// matched text range rather than rewriting the whole paragraph.\n\nconst replacements = [\n  [\"18\u00d743=774\", \"58\u00d726=1508\"],\n  [\"73\u00d718=1314\", \"81\u00d755=4455\"],\n  [\"41\u00d766=2706\", \"68\u00d746=3128\"],\n  [\"94\u00d731=2914\", \"32\u00d785=2720\"],\n  [\"20\u00d737=740\", \"65\u00d780=5200\"],\n  [\"75\u00d794=7050\", \"75\u00d777=5775\"],\n  [\"67\u00d774=4958\", \"79\u00d763=4977\"],\n  [\"17\u00d769=1173\", \"92\u00d735=3220\"],\n  [\"67\u00d740=2680\", \"54\u00d789=4806\"],\n  [\"74\u00d797=7178\", \"25\u00d721=525\"],\n  [\"84\u00d749=4116\", \"93\u00d784=7812\"],\n  [\"23\u00d778=1794\", \"73\u00d769=5037\"],\n  [\"65\u00d762=4030\", \"99\u00d722=2178\"],\n  [\"51\u00d779=4029\", \"48\u00d797=4656\"],\n  [\"34\u00d765=2210\", \"61\u00d762=3782\"],\n  [\"72\u00d728=2016\", \"26\u00d782=2132\"],\n  [\"87\u00d730=2610\", \"55\u00d740=2200\"],\n  [\"62\u00d744=2728\", \"59\u00d756=3304\"],\n  [\"20\u00d795=1900\", \"62\u00d723=1426\"],\n  [\"78\u00d789=6942\", \"82\u00d717=1394\"],\n  [\"81\u00d742=3402\", \"93\u00d747=4371\"],\n  [\"80\u00d741=3280\", \"88\u00d720=1760\"],\n  [\"82\u00d725=2050\", \"66\u00d777=5082\"],\n  [\"47\u00d798=4606\", \"83\u00d760=4980\"],\n  [\"87\u00d779=6873\", \"87\u00d788=7656\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication problems/answers table: each cell's old\n# \"A\u00d7B=C\" text is replaced with a new \"A\u00d7B=C\" text. wdReplaceOne (1)\n# replaces the single occurrence and leaves the surrounding run\n# formatting (font, size, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"18\u00d743=774\";  New = \"58\u00d726=1508\"},\n    @{Old = \"73\u00d718=1314\"; New = \"81\u00d755=4455\"},\n    @{Old = \"41\u00d766=2706\"; New = \"68\u00d746=3128\"},\n    @{Old = \"94\u00d731=2914\"; New = \"32\u00d785=2720\"},\n    @{Old = \"20\u00d737=740\";  New = \"65\u00d780=5200\"},\n    @{Old = \"75\u00d794=7050\"; New = \"75\u00d777=5775\"},\n    @{Old = \"67\u00d774=4958\"; New = \"79\u00d763=4977\"},\n    @{Old = \"17\u00d769=1173\"; New = \"92\u00d735=3220\"},\n    @{Old = \"67\u00d740=2680\"; New = \"54\u00d789=4806\"},\n    @{Old = \"74\u00d797=7178\"; New = \"25\u00d721=525\"},\n    @{Old = \"84\u00d749=4116\"; New = \"93\u00d784=7812\"},\n    @{Old = \"23\u00d778=1794\"; New = \"73\u00d769=5037\"},\n    @{Old = \"65\u00d762=4030\"; New = \"99\u00d722=2178\"},\n    @{Old = \"51\u00d779=4029\"; New = \"48\u00d797=4656\"},\n    @{Old = \"34\u00d765=2210\"; New = \"61\u00d762=3782\"},\n    @{Old = \"72\u00d728=2016\"; New = \"26\u00d782=2132\"},\n    @{Old = \"87\u00d730=2610\"; New = \"55\u00d740=2200\"},\n    @{Old = \"62\u00d744=2728\"; New = \"59\u00d756=3304\"},\n    @{Old = \"20\u00d795=1900\"; New = \"62\u00d723=1426\"},\n    @{Old = \"78\u00d789=6942\"; New = \"82\u00d717=1394\"},\n    @{Old = \"81\u00d742=3402\"; New = \"93\u00d747=4371\"},\n    @{Old = \"80\u00d741=3280\"; New = \"88\u00d720=1760\"},\n    @{Old = \"82\u00d725=2050\"; New = \"66\u00d777=5082\"},\n    @{Old = \"47\u00d798=4606\"; New = \"83\u00d760=4980\"},\n    @{Old = \"87\u00d779=6873\"; New = \"87\u00d788=7656\"}\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 1)\n}\n"}
